# Migrar todo a el codigo a ingles (parcial)
# Adds two new Spanish/English naming-convention rows to the glossary sheet
# ("Arco"/"arc" and "Nodo de servicio"/"broker node"), widens column A to
# fit the longer Spanish terms, and leaves the selection on B7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows 5 and 6 in columns A (Spanish) / B (English)
$ws.Range("A5").Value = "Arco"
$ws.Range("B5").Value = "arc"

$ws.Range("A6").Value = "Nodo de servicio"
$ws.Range("B6").Value = "broker node"

# Widen column A so the new, longer terms are fully visible (stored width 24)
$ws.Columns.Item(1).ColumnWidth = 23.166666666666668

# Move/save the active selection to B7
$ws.Range("B7").Select()
